$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Theme: background (Light1) color change ---
$tcs = $wb.Theme.ThemeColorScheme
$tcs.Colors(2).RGB = 13625548   # CCE8CF

# --- Sheet1 (RabbitMQ): selection change only ---
$ws1.Activate()
$ws1.Range("B25").Select()

# --- Sheet2 (ThreadLocal): fix typo, add row 4, change selection ---
$ws2.Range("B3").Value = "学习ThreadLocal源码，自编写ThreadLocal"
$ws2.Range("A3").Copy()
$ws2.Range("A4").PasteSpecial(-4122)
$ws2.Activate()
$ws2.Range("B24").Select()

# --- Sheet3 (quartz): new sheet with data ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "quartz"

$ws3.Range("A1").Value = 44101
$ws3.Range("B1").Value = "博客学习"
$ws3.Range("C1").Value = "花费3小时"
$ws3.Range("D1").Value = "周四"

$ws3.Range("A2").Value = 44102
$ws3.Range("B2").Value = "源码查看"
$ws3.Range("C2").Value = "花费5小时"
$ws3.Range("D2").Value = "周五"
$ws3.Range("E2").Value = "下载源码，查看使用方式"

$ws3.Range("A3").Value = 44112
$ws3.Range("B3").Value = "quartz 模型搭建"
$ws3.Range("C3").Value = "花费5小时"
$ws3.Range("D3").Value = "周四"

$ws3.Range("A4").Value = 44115
$ws3.Range("B4").Value = "官方文档"
$ws3.Range("C4").Value = "花费3小时"
$ws3.Range("D4").Value = "周日"

$ws3.Range("A5").Value = 44117
$ws3.Range("B5").Value = "运行bug fix"
$ws3.Range("C5").Value = "花费5小时"
$ws3.Range("D5").Value = "周二"
$ws3.Range("E5").Value = "排查问题并验证"

# Reuse the existing date style (numFmtId 14) on column A instead of creating a
# brand-new custom number format entry.
$ws2.Range("A1").Copy()
$ws3.Range("A1:A5").PasteSpecial(-4122)
# Re-enter the values since PasteSpecial(Formats) should not disturb them, but
# make sure dates are still correct after the format-only paste.
$ws3.Range("A1").Value = 44101
$ws3.Range("A2").Value = 44102
$ws3.Range("A3").Value = 44112
$ws3.Range("A4").Value = 44115
$ws3.Range("A5").Value = 44117

# Column widths (bestFit, matching the other sheets' style)
$ws3.Columns.Item(1).ColumnWidth = 11.125
$ws3.Columns.Item(2).ColumnWidth = 15.25
$ws3.Columns.Item(3).ColumnWidth = 10
$ws3.Columns.Item(4).ColumnWidth = 5.25
$ws3.Columns.Item(5).ColumnWidth = 23.5

$ws3.Range("E17").Select()
$ws3.Activate()
